$wb = $excel.ActiveWorkbook

# Worksheets: "USD conversion" (rId1, 1st sheet) and "EUR conversion" (rId2, 2nd sheet)
$wsUsd = $wb.Worksheets.Item("USD conversion")
$wsEur = $wb.Worksheets.Item("EUR conversion")

# --- USD conversion sheet: Zimbabwe row (row 59) now filled in with 1 for every year (B:O) ---
for ($col = 2; $col -le 15; $col++) {
    $wsUsd.Cells.Item(59, $col).Value = 1
}

# --- EUR conversion sheet: fix a precision rounding artifact on Ecuador's 2011 rate (B19) ---
$wsEur.Range("B19").Value = 1.3257166666667

# --- EUR conversion sheet: Zimbabwe row (row 58) copied from Ecuador's row (row 19) ---
$wsEur.Range("B19:O19").Copy()
$wsEur.Range("B58:O58").PasteSpecial(-4163)

# --- Restore view/selection state ---
$wsUsd.Activate()
$wsUsd.Range("L66").Select()

$wsEur.Activate()
$wsEur.Range("B58:O58").Select()
